$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New recomputed s_val data (regenerated after filtering save games).
# Columns: B=TB, C=d2S, D=K, E=IP, G=sum  (F=Win is unchanged)
$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    3  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    4  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    5  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082)
    6  = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 31.61296591696135)
    7  = @(0.127881588408715, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.192886062056065)
    8  = @(0.3048080303191223, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 1.001517020209437)
    9  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    10 = @(0.01514828764759746, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 4.725148573877958)
    11 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    12 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    13 = @(3.230985683306322, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 4.198080582305154)
    14 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    15 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    16 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    17 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
